# Applies the "Updated cryptos list" refresh: per-row Price (D) / Volume(1h) (E)
# updates, plus the ARBITRUM/HuobiToken row swap (rows 50-51, columns B/C).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cellRef -> new value. A leading "'" forces text storage (matches the
# source file, where every data cell is an inline string) so values such as
# "0.0900" or "231.73" keep their exact digits/trailing zeros instead of
# being auto-converted to a Number by Excel's usual Value-assignment parsing.
$updates = [ordered]@{
    "D2" = "41.826.15"
    "E2" = "  +0.44%  "
    "D3" = "2.233.68"
    "E3" = "  -1.05%  "
    "E4" = "  -0.07%  "
    "D5" = "'231.73"
    "E5" = "  -1.09%  "
    "D6" = "'0.623"
    "E6" = "  -2.43%  "
    "D7" = "'60.15"
    "E7" = "  -7.49%  "
    "E8" = "  -0.02%  "
    "D9" = "'0.403"
    "E9" = "  -1.46%  "
    "D10" = "'57.94"
    "E10" = "  -3.32%  "
    "D11" = "'0.0900"
    "E11" = "  -0.77%  "
    "D13" = "2.568.33"
    "E13" = "  -1.11%  "
    "D14" = "'15.50"
    "E14" = "  -4.17%  "
    "D15" = "'22.60"
    "E15" = "  +0.85%  "
    "D16" = "'5.65"
    "E16" = "  -0.40%  "
    "D17" = "'0.801"
    "E17" = "  -4.01%  "
    "D18" = "2.243.00"
    "E18" = "  -1.04%  "
    "D19" = "41.746.58"
    "E19" = "  +0.42%  "
    "D20" = "0.0₃0911"
    "E20" = "  -1.04%  "
    "D21" = "'72.45"
    "E21" = "  -2.19%  "
    "D22" = "'6.14"
    "E22" = "  -0.64%  "
    "D23" = "'247.79"
    "E23" = "  -2.20%  "
    "E24" = "  -0.20%  "
    "D25" = "'2.37"
    "E25" = "  -1.94%  "
    "D26" = "'2.31"
    "E26" = "  -1.28%  "
    "D27" = "'9.76"
    "E27" = "  -0.62%  "
    "D28" = "'169.66"
    "E28" = "  -1.94%  "
    "E29" = "  -2.70%  "
    "D30" = "'19.92"
    "E30" = "  -2.98%  "
    "D31" = "'1.41"
    "E31" = "  -2.16%  "
    "D32" = "'2.58"
    "E32" = "  -8.20%  "
    "D33" = "'0.122"
    "E33" = "  -1.85%  "
    "D34" = "'5.01"
    "E34" = "  +1.56%  "
    "D35" = "'4.70"
    "E35" = "  -1.13%  "
    "D36" = "'0.0656"
    "E36" = "  +3.05%  "
    "D37" = "'6.56"
    "E37" = "  -9.32%  "
    "D38" = "'2.41"
    "E38" = "  -1.95%  "
    "D39" = "'3.61"
    "E39" = "  -6.86%  "
    "D40" = "'0.000243"
    "E40" = "  +3.68%  "
    "E41" = "  +0.09%  "
    "E42" = "  +0.89%  "
    "D43" = "'8.67"
    "E43" = "  -1.27%  "
    "E44" = "  -1.10%  "
    "D45" = "'99.19"
    "D46" = "'0.0966"
    "E46" = "  +2.48%  "
    "D47" = "'4.42"
    "E47" = "  -9.44%  "
    "D48" = "1.475.12"
    "E48" = "  -2.66%  "
    "D49" = "'16.67"
    "E49" = "  -8.92%  "
    "B50" = "HuobiToken"
    "C50" = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
    "D50" = "'2.77"
    "E50" = "  -1.48%  "
    "B51" = "ARBITRUM"
    "C51" = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
    "D51" = "'1.08"
    "E51" = "  -2.84%  "
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}

Write-Output "Applied $($updates.Count) cell updates"
